# Fix duong dan python cho Windows
# This edit adds a new user (68) to the cosine similarity matrix:
#   - a new column AO (header 68, similarity values, self = 1 at AO41)
#   - a new row 41 (header 68, similarity values, self = 1 at AO41)
# and refreshes the previously-computed row/column for user 33 (row 9 / column I),
# whose similarity values changed now that an extra user participates in the
# calculation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- 1. Grow the used range: new header cells for user 68 ----
$ws.Range("AO1").Value2 = 68
$ws.Range("A41").Value2 = 68

# Match the bold / centered / thin-bordered header style used by the other
# header cells (row 1 and column A).
foreach ($rng in @($ws.Range("AO1"), $ws.Range("A41"))) {
    $rng.Font.Bold = $true
    $rng.HorizontalAlignment = -4108   # xlCenter
    $rng.VerticalAlignment = -4160     # xlTop
    $rng.Borders.LineStyle = 1         # xlContinuous (thin box border)
}

# ---- 2. Recomputed similarities for user 33 (row 9 / column I) ----
$row9 = New-Object 'object[,]' 1,40
$row9Vals = @(0.01,0.04,0.08,0.23,0.54,0,0.01,1,0.01,0.09,0,0.01,0.01,0,0,0,0,0.31,0,0.05,0.01,0,0,0,0,0.29,0.11,0.44,0.02,0.11,0.13,0,0.16,0.08,0.11,0,0.04,0,0.01,0)
for ($i = 0; $i -lt $row9Vals.Length; $i++) { $row9[0,$i] = $row9Vals[$i] }
$ws.Range("B9:AO9").Value2 = $row9

$colI = New-Object 'object[,]' 40,1
$colIVals = @(0.01,0.04,0.08,0.23,0.54,0,0.01,1,0.01,0.09,0,0.01,0.01,0,0,0,0,0.31,0,0.05,0.01,0,0,0,0,0.29,0.11,0.44,0.02,0.11,0.13,0,0.16,0.08,0.11,0,0.04,0,0.01,0)
for ($i = 0; $i -lt $colIVals.Length; $i++) { $colI[$i,0] = $colIVals[$i] }
$ws.Range("I2:I41").Value2 = $colI

# ---- 3. New column AO (user 68) similarities for existing users (rows 2-40) ----
$aoCol = New-Object 'object[,]' 39,1
$aoVals = @(0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)
for ($i = 0; $i -lt $aoVals.Length; $i++) { $aoCol[$i,0] = $aoVals[$i] }
$ws.Range("AO2:AO40").Value2 = $aoCol

# ---- 4. New row 41 (user 68) similarities for existing users (columns B-AN) ----
$row41 = New-Object 'object[,]' 1,39
$row41Vals = @(0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)
for ($i = 0; $i -lt $row41Vals.Length; $i++) { $row41[0,$i] = $row41Vals[$i] }
$ws.Range("B41:AN41").Value2 = $row41

# ---- 5. Self-similarity of the new user 68 ----
$ws.Range("AO41").Value2 = 1

Write-Output "applied cosine similarity matrix update"
